$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change fund to right: rename FUNDCODE column header to RIGHTID
$ws.Range("D1").Value = "RIGHTID"

# Update the related validation/description text for the RIGHTID column
$ws.Range("E3").Value = "RIGHTID จะต้องมีในฐานข้อมูล"
